$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"): copy the formatting of the
# existing header cell H1 (bold font, thin border, centered alignment)
# so the new header cells share the same style as the other headers.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells I2 and J2 (plain numeric values, same as H2)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
